$wb = $excel.ActiveWorkbook

# This script applies updated price/profit figures to the Leve profit
# tracking sheets, as produced by the scheduled data-refresh runner.

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 345
$ws.Range("I2").Value = 400
$ws.Range("K2").Value = 400
$ws.Range("M2").Value = -287
$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 1625
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 1625
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -1450
$ws.Range("N40").Value = -2850
$ws.Range("H64").Value = 3767.96
$ws.Range("I64").Value = 3630.7693
$ws.Range("K64").Value = 3630.7693
$ws.Range("M64").Value = -3382.7693
$ws.Range("H67").Value = 3767.96
$ws.Range("I67").Value = 3630.7693
$ws.Range("K67").Value = 3630.7693
$ws.Range("M67").Value = -2772.7693
$ws.Range("H106").Value = 4355.077
$ws.Range("I106").Value = 2885.8333
$ws.Range("J106").Value = 5614.4287
$ws.Range("K106").Value = 2885.8333
$ws.Range("L106").Value = 5614.4287
$ws.Range("M106").Value = -2254.8333
$ws.Range("N106").Value = -6876.4287
$ws.Range("H137").Value = 2040.7632
$ws.Range("I137").Value = 1022.12
$ws.Range("J137").Value = 3999.6924
$ws.Range("K137").Value = 3066.36
$ws.Range("L137").Value = 11999.0772
$ws.Range("M137").Value = -516.3600000000001
$ws.Range("N137").Value = -17099.0772

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13336942
$ws.Range("I32").Value = 14927722
$ws.Range("J32").Value = 14151.5
$ws.Range("K32").Value = 14927722
$ws.Range("L32").Value = 14151.5
$ws.Range("M32").Value = -14927435
$ws.Range("N32").Value = -14725.5
$ws.Range("H44").Value = 24000
$ws.Range("J44").Value = 24000
$ws.Range("L44").Value = 24000
$ws.Range("N44").Value = -24976
$ws.Range("H54").Value = 16024.5
$ws.Range("J54").Value = 16024.5
$ws.Range("L54").Value = 16024.5
$ws.Range("N54").Value = -17562.5
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = $null
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = $null
$ws.Range("H122").Value = 1913.5
$ws.Range("I122").Value = 1795
$ws.Range("J122").Value = 2032
$ws.Range("K122").Value = 5385
$ws.Range("L122").Value = 6096
$ws.Range("M122").Value = -2935
$ws.Range("N122").Value = -10996
$ws.Range("H132").Value = 1502.725
$ws.Range("I132").Value = 1138.8064
$ws.Range("K132").Value = 3416.4192
$ws.Range("M132").Value = -886.4191999999998

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1376.5
$ws.Range("I86").Value = 1499.6666
$ws.Range("J86").Value = 1007
$ws.Range("K86").Value = 1499.6666
$ws.Range("L86").Value = 1007
$ws.Range("M86").Value = -376.6666
$ws.Range("N86").Value = -3253
$ws.Range("H89").Value = 1376.5
$ws.Range("I89").Value = 1499.6666
$ws.Range("J89").Value = 1007
$ws.Range("K89").Value = 7498.333000000001
$ws.Range("L89").Value = 5035
$ws.Range("M89").Value = -1882.333000000001
$ws.Range("N89").Value = -16267
$ws.Range("H105").Value = 2513.66
$ws.Range("I105").Value = 1317.762
$ws.Range("J105").Value = 2831.557
$ws.Range("K105").Value = 1317.762
$ws.Range("L105").Value = 2831.557
$ws.Range("M105").Value = 429.2380000000001
$ws.Range("N105").Value = -6325.557
$ws.Range("H107").Value = 1892.1471
$ws.Range("I107").Value = 2175.762
$ws.Range("J107").Value = 1434
$ws.Range("K107").Value = 2175.762
$ws.Range("L107").Value = 1434
$ws.Range("M107").Value = -255.7620000000002
$ws.Range("N107").Value = -5274
$ws.Range("H129").Value = 35499.5
$ws.Range("J129").Value = 35499.5
$ws.Range("L129").Value = 35499.5
$ws.Range("N129").Value = -45499.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1000
$ws.Range("I105").Value = 1400
$ws.Range("J105").Value = 800
$ws.Range("K105").Value = 1400
$ws.Range("L105").Value = 800
$ws.Range("M105").Value = 347
$ws.Range("N105").Value = -4294
$ws.Range("H107").Value = 615.75
$ws.Range("I107").Value = 456.9091
$ws.Range("J107").Value = 809.8889
$ws.Range("K107").Value = 456.9091
$ws.Range("L107").Value = 809.8889
$ws.Range("M107").Value = 1463.0909
$ws.Range("N107").Value = -4649.8889

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 15.55
$ws.Range("I12").Value = 33.714287
$ws.Range("K12").Value = 101.142861
$ws.Range("M12").Value = 71.857139
$ws.Range("H69").Value = 2603.6875
$ws.Range("I69").Value = 2499.5
$ws.Range("J69").Value = 2618.5715
$ws.Range("K69").Value = 7498.5
$ws.Range("L69").Value = 7855.7145
$ws.Range("M69").Value = -6687.5
$ws.Range("N69").Value = -9477.7145
$ws.Range("H72").Value = 2603.6875
$ws.Range("I72").Value = 2499.5
$ws.Range("J72").Value = 2618.5715
$ws.Range("K72").Value = 22495.5
$ws.Range("L72").Value = 23567.1435
$ws.Range("M72").Value = -18439.5
$ws.Range("N72").Value = -31679.1435
$ws.Range("H113").Value = 766895.4
$ws.Range("I113").Value = 1379896.9
$ws.Range("K113").Value = 4139690.7
$ws.Range("M113").Value = -4137520.7
$ws.Range("H129").Value = 6214.478
$ws.Range("I129").Value = 1817.8667
$ws.Range("J129").Value = 14458.125
$ws.Range("K129").Value = 5453.6001
$ws.Range("L129").Value = 43374.375
$ws.Range("M129").Value = -453.6000999999997
$ws.Range("N129").Value = -53374.375
$ws.Range("H132").Value = 778569.75
$ws.Range("J132").Value = 1685154.6
$ws.Range("L132").Value = 15166391.4
$ws.Range("N132").Value = -15171451.4

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6789.9565
$ws.Range("I70").Value = 7325.5557
$ws.Range("J70").Value = 4861.8
$ws.Range("K70").Value = 7325.5557
$ws.Range("L70").Value = 4861.8
$ws.Range("M70").Value = -7055.5557
$ws.Range("N70").Value = -5401.8
$ws.Range("H73").Value = 6789.9565
$ws.Range("I73").Value = 7325.5557
$ws.Range("J73").Value = 4861.8
$ws.Range("K73").Value = 7325.5557
$ws.Range("L73").Value = 4861.8
$ws.Range("M73").Value = -6389.5557
$ws.Range("N73").Value = -6733.8
$ws.Range("H109").Value = 10285
$ws.Range("J109").Value = 10285
$ws.Range("L109").Value = 10285
$ws.Range("N109").Value = -12365
$ws.Range("H122").Value = 2730
$ws.Range("I122").Value = 2595
$ws.Range("K122").Value = 7785
$ws.Range("M122").Value = -5335

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 640
$ws.Range("J22").Value = 600
$ws.Range("L22").Value = 600
$ws.Range("N22").Value = -1190
$ws.Range("H27").Value = 640
$ws.Range("J27").Value = 600
$ws.Range("L27").Value = 600
$ws.Range("N27").Value = -814
$ws.Range("H46").Value = 1125
$ws.Range("I46").Value = 937.5
$ws.Range("J46").Value = 1500
$ws.Range("K46").Value = 937.5
$ws.Range("L46").Value = 1500
$ws.Range("M46").Value = -749.5
$ws.Range("N46").Value = -1876

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1158.8368
$ws.Range("I132").Value = 728.4211
$ws.Range("J132").Value = 2645.7273
$ws.Range("K132").Value = 2185.2633
$ws.Range("L132").Value = 7937.1819
$ws.Range("M132").Value = 344.7366999999999
$ws.Range("N132").Value = -12997.1819
